$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.16594766666667
$ws.Range("H2").Value = 39.497843
$ws.Range("I2").Value = 0.6940777873489595
$ws.Range("J2").Value = 0.6940777873489595
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3407069999999999
$ws.Range("N2").Value = 1.022121
$ws.Range("O2").Value = 0.1055965976712818
$ws.Range("P2").Value = 0.1055965976712818
$ws.Range("Q2").Value = 4.485730531667
$ws.Range("R2").Value = 40.371574785003
$ws.Range("S2").Value = 0.07329225286326156
$ws.Range("T2").Value = 0.07329225286326155
$ws.Range("G3").Value = 13.16594766666667
$ws.Range("H3").Value = 39.497843
$ws.Range("I3").Value = 0.6940777873489595
$ws.Range("J3").Value = 0.6940777873489595
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.220310333333333
$ws.Range("N3").Value = 3.660931
$ws.Range("O3").Value = 0.3782153560188308
$ws.Range("P3").Value = 0.3782153560188308
$ws.Range("Q3").Value = 16.06654198575922
$ws.Range("R3").Value = 144.598877871833
$ws.Range("S3").Value = 0.2625108774469491
$ws.Range("T3").Value = 0.2625108774469491
$ws.Range("G4").Value = 13.16594766666667
$ws.Range("H4").Value = 39.497843
$ws.Range("I4").Value = 0.6940777873489595
$ws.Range("J4").Value = 0.6940777873489595
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.665478666666666
$ws.Range("N4").Value = 4.996435999999999
$ws.Range("O4").Value = 0.5161880463098875
$ws.Range("P4").Value = 0.5161880463098875
$ws.Range("Q4").Value = 21.92760496528311
$ws.Range("R4").Value = 197.348444687548
$ws.Range("S4").Value = 0.3582746570387489
$ws.Range("T4").Value = 0.3582746570387489
$ws.Range("G5").Value = 3.23724
$ws.Range("H5").Value = 9.71172
$ws.Range("I5").Value = 0.1706596770095176
$ws.Range("J5").Value = 0.1706596770095176
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3407069999999999
$ws.Range("N5").Value = 1.022121
$ws.Range("O5").Value = 0.1055965976712818
$ws.Range("P5").Value = 0.1055965976712818
$ws.Range("Q5").Value = 1.10295032868
$ws.Range("R5").Value = 9.926552958119998
$ws.Range("S5").Value = 0.01802108125188493
$ws.Range("T5").Value = 0.01802108125188493
$ws.Range("G6").Value = 3.23724
$ws.Range("H6").Value = 9.71172
$ws.Range("I6").Value = 0.1706596770095176
$ws.Range("J6").Value = 0.1706596770095176
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.220310333333333
$ws.Range("N6").Value = 3.660931
$ws.Range("O6").Value = 0.3782153560188308
$ws.Range("P6").Value = 0.3782153560188308
$ws.Range("Q6").Value = 3.950437423479999
$ws.Range("R6").Value = 35.55393681131999
$ws.Range("S6").Value = 0.06454611049821338
$ws.Range("T6").Value = 0.06454611049821338
$ws.Range("G7").Value = 3.23724
$ws.Range("H7").Value = 9.71172
$ws.Range("I7").Value = 0.1706596770095176
$ws.Range("J7").Value = 0.1706596770095176
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.665478666666666
$ws.Range("N7").Value = 4.996435999999999
$ws.Range("O7").Value = 0.5161880463098875
$ws.Range("P7").Value = 0.5161880463098875
$ws.Range("Q7").Value = 5.391554158879999
$ws.Range("R7").Value = 48.52398742991999
$ws.Range("S7").Value = 0.08809248525941932
$ws.Range("T7").Value = 0.08809248525941932
$ws.Range("G8").Value = 2.565792333333333
$ws.Range("H8").Value = 7.697376999999999
$ws.Range("I8").Value = 0.1352625356415228
$ws.Range("J8").Value = 0.1352625356415228
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3407069999999999
$ws.Range("N8").Value = 1.022121
$ws.Range("O8").Value = 0.1055965976712818
$ws.Range("P8").Value = 0.1055965976712818
$ws.Range("Q8").Value = 0.8741834085129997
$ws.Range("R8").Value = 7.867650676616998
$ws.Range("S8").Value = 0.01428326355613529
$ws.Range("T8").Value = 0.01428326355613529
$ws.Range("G9").Value = 2.565792333333333
$ws.Range("H9").Value = 7.697376999999999
$ws.Range("I9").Value = 0.1352625356415228
$ws.Range("J9").Value = 0.1352625356415228
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.220310333333333
$ws.Range("N9").Value = 3.660931
$ws.Range("O9").Value = 0.3782153560188308
$ws.Range("P9").Value = 0.3782153560188308
$ws.Range("Q9").Value = 3.13106289755411
$ws.Range("R9").Value = 28.17956607798699
$ws.Range("S9").Value = 0.05115836807366832
$ws.Range("T9").Value = 0.05115836807366833
$ws.Range("G10").Value = 2.565792333333333
$ws.Range("H10").Value = 7.697376999999999
$ws.Range("I10").Value = 0.1352625356415228
$ws.Range("J10").Value = 0.1352625356415228
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.665478666666666
$ws.Range("N10").Value = 4.996435999999999
$ws.Range("O10").Value = 0.5161880463098875
$ws.Range("P10").Value = 0.5161880463098875
$ws.Range("Q10").Value = 4.273272394263554
$ws.Range("R10").Value = 38.45945154837199
$ws.Range("S10").Value = 0.06982090401171916
$ws.Range("T10").Value = 0.06982090401171917
